$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: record an end time (C16) for the existing start time (B16),
# and log the entry text in column E.
$ws.Range("C16").Value = 0.80208333333333337
$ws.Range("C16").NumberFormat = $ws.Range("B16").NumberFormat
$ws.Range("E16").Value = "attempt more flext, sudden idea for min and backpivot"

# Row 17: new start/end time pair plus its log entry.
$ws.Range("B17").Value = 0.83333333333333337
$ws.Range("B17").NumberFormat = $ws.Range("B16").NumberFormat
$ws.Range("C17").Value = 0.88888888888888884
$ws.Range("C17").NumberFormat = $ws.Range("B16").NumberFormat
$ws.Range("E17").Value = "Woohoo! Now able to build resonance object files from inside min devkit! Now work on buffer conversion"

# New subtotal for this block of entries (rows 13-16) in E25,
# styled like the other wrapped time-format cells.
$ws.Range("E25").Formula = "=SUM(D13:D16)"

# Move the active selection.
$ws.Range("E18").Select()
